$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Copy formatting (styles) from row 410 into new rows 411-420
$ws.Range("A410:M410").Copy()
$ws.Range("A411:M420").PasteSpecial(-4122)

# 2) Force text columns to Text format temporarily so numeric-looking
#    strings (e.g. "2006", "9010") are stored as shared strings, not numbers
$ws.Range("A411:A420").NumberFormat = "@"
$ws.Range("D411:D420").NumberFormat = "@"
$ws.Range("E411:E420").NumberFormat = "@"
$ws.Range("F411:F420").NumberFormat = "@"
$ws.Range("G411:G420").NumberFormat = "@"
$ws.Range("I411:I420").NumberFormat = "@"
$ws.Range("L411:L420").NumberFormat = "@"
$ws.Range("M411:M420").NumberFormat = "@"

# 3) Populate cell values row by row (order matters for shared-string table order)
# Row 411
$ws.Range("A411").Value = "●"
$ws.Range("D411").Value = "9010"
$ws.Range("E411").Value = "Bacteria:Binomial (genus species)"
$ws.Range("F411").Value = "1: 2617"
$ws.Range("G411").Value = "1: 2637"
$ws.Range("H411").Value = 0
$ws.Range("I411").Value = "Klebsiella pneumoniae"
$ws.Range("J411").Value = 21
$ws.Range("K411").Value = 0.083287
$ws.Range("L411").Value = "Sonia"
$ws.Range("M411").Value = "11/8/18 14:35:00"

# Row 412
$ws.Range("A412").Value = "●"
$ws.Range("D412").Value = "21726"
$ws.Range("E412").Value = "Bacteria:Binomial (genus species)"
$ws.Range("F412").Value = "1: 2873"
$ws.Range("G412").Value = "1: 2900"
$ws.Range("H412").Value = 0
$ws.Range("I412").Value = "ycobacterium tubercu- `nlosis"
$ws.Range("J412").Value = 28
$ws.Range("K412").Value = 0.094362
$ws.Range("L412").Value = "Sonia"
$ws.Range("M412").Value = "11/8/18 14:35:00"

# Row 413
$ws.Range("A413").Value = "●"
$ws.Range("D413").Value = "21726"
$ws.Range("E413").Value = "Event month"
$ws.Range("F413").Value = "3: 2050"
$ws.Range("G413").Value = "3: 2056"
$ws.Range("H413").Value = 0
$ws.Range("I413").Value = "October"
$ws.Range("J413").Value = 7
$ws.Range("K413").Value = 0.02359
$ws.Range("L413").Value = "Sonia"
$ws.Range("M413").Value = "11/12/18 14:04:00"

# Row 414
$ws.Range("A414").Value = "●"
$ws.Range("D414").Value = "21726"
$ws.Range("E414").Value = "Event month"
$ws.Range("F414").Value = "3: 2066"
$ws.Range("G414").Value = "3: 2072"
$ws.Range("H414").Value = 0
$ws.Range("I414").Value = "October"
$ws.Range("J414").Value = 7
$ws.Range("K414").Value = 0.02359
$ws.Range("L414").Value = "Sonia"
$ws.Range("M414").Value = "11/12/18 14:04:00"

# Row 415
$ws.Range("A415").Value = "●"
$ws.Range("D415").Value = "21726"
$ws.Range("E415").Value = "Event year"
$ws.Range("F415").Value = "3: 2058"
$ws.Range("G415").Value = "3: 2061"
$ws.Range("H415").Value = 0
$ws.Range("I415").Value = "2006"
$ws.Range("J415").Value = 4
$ws.Range("K415").Value = 0.01348
$ws.Range("L415").Value = "Sonia"
$ws.Range("M415").Value = "11/12/18 14:04:00"

# Row 416
$ws.Range("A416").Value = "●"
$ws.Range("D416").Value = "21726"
$ws.Range("E416").Value = "Event year"
$ws.Range("F416").Value = "3: 2074"
$ws.Range("G416").Value = "3: 2077"
$ws.Range("H416").Value = 0
$ws.Range("I416").Value = "2008"
$ws.Range("J416").Value = 4
$ws.Range("K416").Value = 0.01348
$ws.Range("L416").Value = "Sonia"
$ws.Range("M416").Value = "11/12/18 14:04:00"

# Row 417
$ws.Range("A417").Value = "●"
$ws.Range("D417").Value = "21726"
$ws.Range("E417").Value = "B"
$ws.Range("F417").Value = "3: 2074"
$ws.Range("G417").Value = "3: 2077"
$ws.Range("H417").Value = 0
$ws.Range("I417").Value = "2008"
$ws.Range("J417").Value = 4
$ws.Range("K417").Value = 0.01348
$ws.Range("L417").Value = "Sonia"
$ws.Range("M417").Value = "11/12/18 14:05:00"

# Row 418
$ws.Range("A418").Value = "●"
$ws.Range("D418").Value = "21726"
$ws.Range("E418").Value = "B"
$ws.Range("F418").Value = "3: 2066"
$ws.Range("G418").Value = "3: 2072"
$ws.Range("H418").Value = 0
$ws.Range("I418").Value = "October"
$ws.Range("J418").Value = 7
$ws.Range("K418").Value = 0.02359
$ws.Range("L418").Value = "Sonia"
$ws.Range("M418").Value = "11/12/18 14:05:00"

# Row 419
$ws.Range("A419").Value = "●"
$ws.Range("D419").Value = "21726"
$ws.Range("E419").Value = "A"
$ws.Range("F419").Value = "3: 2058"
$ws.Range("G419").Value = "3: 2061"
$ws.Range("H419").Value = 0
$ws.Range("I419").Value = "2006"
$ws.Range("J419").Value = 4
$ws.Range("K419").Value = 0.01348
$ws.Range("L419").Value = "Sonia"
$ws.Range("M419").Value = "11/12/18 14:05:00"

# Row 420
$ws.Range("A420").Value = "●"
$ws.Range("D420").Value = "21726"
$ws.Range("E420").Value = "A"
$ws.Range("F420").Value = "3: 2050"
$ws.Range("G420").Value = "3: 2056"
$ws.Range("H420").Value = 0
$ws.Range("I420").Value = "October"
$ws.Range("J420").Value = 7
$ws.Range("K420").Value = 0.02359
$ws.Range("L420").Value = "Sonia"
$ws.Range("M420").Value = "11/12/18 14:05:00"

# 4) Restore exact original cell styles (overwritten by the NumberFormat tweak above)
$ws.Range("A410:M410").Copy()
$ws.Range("A411:M420").PasteSpecial(-4122)

# 5) Set explicit row heights to match target layout
$ws.Rows.Item(411).RowHeight = 16
$ws.Rows.Item(412).RowHeight = 30
$ws.Rows.Item(413).RowHeight = 16
$ws.Rows.Item(414).RowHeight = 16
$ws.Rows.Item(415).RowHeight = 16
$ws.Rows.Item(416).RowHeight = 16
$ws.Rows.Item(417).RowHeight = 16
$ws.Rows.Item(418).RowHeight = 16
$ws.Rows.Item(419).RowHeight = 16
$ws.Rows.Item(420).RowHeight = 16

$excel.CutCopyMode = $false
